# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the "Start Ruby" / "Pomelo" series
# (row 107), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(107).Insert()

$ws.Range("A107").Value = 4
$ws.Range("B107").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C107").Value = "Los Lagos"
$ws.Range("D107").Value = 44468
$ws.Range("E107").Value = 10
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100102
$ws.Range("H107").Value = "Cítricos"
$ws.Range("I107").Value = 100102006
$ws.Range("J107").Value = "Pomelo"
$ws.Range("K107").Value = "Start Ruby"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 50
$ws.Range("N107").Value = 12000
$ws.Range("O107").Value = 12000
$ws.Range("P107").Value = 12000
$ws.Range("Q107").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R107").Value = "Región de O'Higgins"
$ws.Range("S107").Value = 857
$ws.Range("T107").Value = 14
